# Update mods data [2026-02-13 15:30:30]
# Append a new row (row 95) to the ModCounts sheet with today's data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 95
$lastRow = $newRow - 1

# Set the new values. The date is written with a leading apostrophe so it is
# stored as literal text ("2026/02/13") instead of being auto-converted into
# a date serial number, matching how the existing date column is stored.
$ws.Cells.Item($newRow, 1).Value = "'2026/02/13"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1200

# Copy the formatting (cell style) from the row above onto the new row so it
# matches the rest of the data rows (centered alignment), without disturbing
# the text values we just set.
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)
